# Weekly update: insert the new week's price row for
# "Vega Modelo de Temuco - Rabanito" ahead of the existing history,
# pushing the older rows (105-113) down to (106-114).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 105; rows 105-113 shift down to 106-114.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the latest week's data.
$ws.Cells.Item(105, 1).Value = 10
$ws.Cells.Item(105, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(105, 3).Value = "La Araucanía"
$ws.Cells.Item(105, 4).Value = 45106
$ws.Cells.Item(105, 5).Value = 9
$ws.Cells.Item(105, 6).Value = 300000001
$ws.Cells.Item(105, 7).Value = "Rabanito"
$ws.Cells.Item(105, 8).Value = "Sin especificar"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 65
$ws.Cells.Item(105, 11).Value = 7000
$ws.Cells.Item(105, 12).Value = 7000
$ws.Cells.Item(105, 13).Value = 7000
$ws.Cells.Item(105, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(105, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(105, 16).Value = 583
$ws.Cells.Item(105, 17).Value = 12
$ws.Cells.Item(105, 18).Value = "Hortaliza"
